$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.026932387521171
$ws.Range("D2").Value = 1.036027381065861
$ws.Range("E2").Value = 1.030562713039575
$ws.Range("F2").Value = 1.042797098461808
$ws.Range("I2").Value = 1.032929471634726
$ws.Range("J2").Value = 1.032093075493366
$ws.Range("K2").Value = 1.038822297994701
$ws.Range("L2").Value = 1.033373363740194
$ws.Range("M2").Value = 1.04557277237106
$ws.Range("N2").Value = 1.01471034057442
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.027914111079139
$ws.Range("D3").Value = 1.03680922309728
$ws.Range("E3").Value = 1.031492008187135
$ws.Range("F3").Value = 1.043766676174212
$ws.Range("I3").Value = 1.033123224824504
$ws.Range("J3").Value = 1.032714603400047
$ws.Range("K3").Value = 1.039413730195474
$ws.Range("L3").Value = 1.034110715006926
$ws.Range("M3").Value = 1.046352839301885
$ws.Range("N3").Value = 1.014917314653295
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.028549518893934
$ws.Range("D4").Value = 1.037314830967388
$ws.Range("E4").Value = 1.032093848025518
$ws.Range("F4").Value = 1.044394290338766
$ws.Range("I4").Value = 1.033246702799699
$ws.Range("J4").Value = 1.033116358741358
$ws.Range("K4").Value = 1.039795475133688
$ws.Range("L4").Value = 1.034587731283177
$ws.Range("M4").Value = 1.046857210457465
$ws.Range("N4").Value = 1.015051056158056
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.028816683828749
$ws.Range("D5").Value = 1.03752731649736
$ws.Range("E5").Value = 1.032346985963186
$ws.Range("F5").Value = 1.044658194016715
$ws.Range("I5").Value = 1.033298159083599
$ws.Range("J5").Value = 1.033285156706183
$ws.Range("K5").Value = 1.039955732111675
$ws.Range("L5").Value = 1.034788244378481
$ws.Range("M5").Value = 1.047069155176819
$ws.Range("N5").Value = 1.015107236629846
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.028861544298123
$ws.Range("D6").Value = 1.037562989497314
$ws.Range("E6").Value = 1.032389496209108
$ws.Range("F6").Value = 1.044702507783259
$ws.Range("I6").Value = 1.033306772188205
$ws.Range("J6").Value = 1.033313492744399
$ws.Range("K6").Value = 1.039982626533336
$ws.Range("L6").Value = 1.034821909947434
$ws.Range("M6").Value = 1.047104736152709
$ws.Range("N6").Value = 1.015116666963372
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.028553088610324
$ws.Range("D7").Value = 1.037317670493828
$ws.Range("E7").Value = 1.032097229979103
$ws.Range("F7").Value = 1.044397816419012
$ws.Range("I7").Value = 1.033247392145191
$ws.Range("J7").Value = 1.033118614621746
$ws.Range("K7").Value = 1.039797617394589
$ws.Range("L7").Value = 1.034590410646404
$ws.Range("M7").Value = 1.046860042839423
$ws.Range("N7").Value = 1.015051807019298
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.02726413071148
$ws.Range("D8").Value = 1.036291668908758
$ws.Range("E8").Value = 1.030876663693603
$ws.Range("F8").Value = 1.04312472304423
$ws.Range("I8").Value = 1.032995343093311
$ws.Range("J8").Value = 1.032303209231828
$ws.Range("K8").Value = 1.0390223717194
$ws.Range("L8").Value = 1.033622575148606
$ws.Range("M8").Value = 1.045836478670991
$ws.Range("N8").Value = 1.01478032648567
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.024994122909269
$ws.Range("D9").Value = 1.034481502763193
$ws.Range("E9").Value = 1.0287299210981
$ws.Range("F9").Value = 1.040883194084553
$ws.Range("I9").Value = 1.032536726116828
$ws.Range("J9").Value = 1.030863213813366
$ws.Range("K9").Value = 1.037649042044032
$ws.Range("L9").Value = 1.031916398635029
$ws.Range("M9").Value = 1.044029924110577
$ws.Range("N9").Value = 1.014300541423273
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.023481685309714
$ws.Range("D10").Value = 1.033273295834974
$ws.Range("E10").Value = 1.027301535061857
$ws.Range("F10").Value = 1.039390123940946
$ws.Range("I10").Value = 1.032221283734639
$ws.Range("J10").Value = 1.029901139512996
$ws.Range("K10").Value = 1.036728664973952
$ws.Range("L10").Value = 1.030778502850209
$ws.Range("M10").Value = 1.042823653569831
$ws.Range("N10").Value = 1.013979757786975
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.022827001572029
$ws.Range("D11").Value = 1.032749802598758
$ws.Range("E11").Value = 1.026683696838485
$ws.Range("F11").Value = 1.038743923382244
$ws.Range("I11").Value = 1.032082400335844
$ws.Range("J11").Value = 1.029484065588566
$ws.Range("K11").Value = 1.036328997145469
$ws.Range("L11").Value = 1.030285684600244
$ws.Range("M11").Value = 1.042300884474629
$ws.Range("N11").Value = 1.013840638173576
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.022583854976826
$ws.Range("D12").Value = 1.032555305033998
$ws.Range("E12").Value = 1.026454304492041
$ws.Range("F12").Value = 1.03850394291285
$ws.Range("I12").Value = 1.032030468584717
$ws.Range("J12").Value = 1.029329072782418
$ws.Range("K12").Value = 1.036180372390807
$ws.Range("L12").Value = 1.030102615288868
$ws.Range("M12").Value = 1.042106638498541
$ws.Range("N12").Value = 1.013788930403985
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.022636009289941
$ws.Range("D13").Value = 1.032597027617922
$ws.Range("E13").Value = 1.026503505378924
$ws.Range("F13").Value = 1.038555417389653
$ws.Range("I13").Value = 1.032041623692887
$ws.Range("J13").Value = 1.029362322559242
$ws.Range("K13").Value = 1.036212260601458
$ws.Range("L13").Value = 1.030141884931339
$ws.Range("M13").Value = 1.042148307906158
$ws.Range("N13").Value = 1.013800023365428
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.022806902342245
$ws.Range("D14").Value = 1.03273372636075
$ws.Range("E14").Value = 1.026664733145482
$ws.Range("F14").Value = 1.038724085554412
$ws.Range("I14").Value = 1.032078114661346
$ws.Range("J14").Value = 1.029471255320485
$ws.Range("K14").Value = 1.036316715247401
$ws.Range("L14").Value = 1.030270552329481
$ws.Range("M14").Value = 1.042284829388987
$ws.Range("N14").Value = 1.013836364655922
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.022912199453363
$ws.Range("D15").Value = 1.03281794452092
$ws.Range("E15").Value = 1.026764084204482
$ws.Range("F15").Value = 1.038828013854973
$ws.Range("I15").Value = 1.032100552346416
$ws.Range("J15").Value = 1.029538362735617
$ws.Range("K15").Value = 1.036381050653428
$ws.Range("L15").Value = 1.030349826634936
$ws.Range("M15").Value = 1.042368936022795
$ws.Range("N15").Value = 1.013858751417916
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.023525138998404
$ws.Range("D16").Value = 1.033308031451348
$ws.Range("E16").Value = 1.027342552970609
$ws.Range("F16").Value = 1.039433016716873
$ws.Range("I16").Value = 1.032230452662234
$ws.Range("J16").Value = 1.029928809071353
$ws.Range("K16").Value = 1.036755165686191
$ws.Range("L16").Value = 1.030811207482362
$ws.Range("M16").Value = 1.042858338700366
$ws.Range("N16").Value = 1.013988986120227
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.02390967666779
$ws.Range("D17").Value = 1.033615361919477
$ws.Range("E17").Value = 1.027705589019194
$ws.Range("F17").Value = 1.039812602245058
$ws.Range("I17").Value = 1.032311321671642
$ws.Range("J17").Value = 1.030173595102369
$ws.Range("K17").Value = 1.036989533865526
$ws.Range("L17").Value = 1.031100592674471
$ws.Range("M17").Value = 1.043165209195013
$ws.Range("N17").Value = 1.014070620656467
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.024133991393626
$ws.Range("D18").Value = 1.033794590517213
$ws.Range("E18").Value = 1.027917405660768
$ws.Range("F18").Value = 1.040034037788921
$ws.Range("I18").Value = 1.032358269752453
$ws.Range("J18").Value = 1.030316327395975
$ws.Range("K18").Value = 1.037126126964476
$ws.Range("L18").Value = 1.031269376290123
$ws.Range("M18").Value = 1.043344158481236
$ws.Range("N18").Value = 1.014118215633658
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.024210480371363
$ws.Range("D19").Value = 1.033855697369304
$ws.Range("E19").Value = 1.027989640469483
$ws.Range("F19").Value = 1.040109546604325
$ws.Range("I19").Value = 1.032374240259023
$ws.Range("J19").Value = 1.03036498738729
$ws.Range("K19").Value = 1.03717268301938
$ws.Range("L19").Value = 1.031326925447695
$ws.Range("M19").Value = 1.043405168233198
$ws.Range("N19").Value = 1.014134440710674
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.023868417301911
$ws.Range("D20").Value = 1.033582391578377
$ws.Range("E20").Value = 1.027666632076317
$ws.Range("F20").Value = 1.039771873229413
$ws.Range("I20").Value = 1.032302668093693
$ws.Range("J20").Value = 1.030147336761847
$ws.Range("K20").Value = 1.03696439973961
$ws.Range("L20").Value = 1.031069545398414
$ws.Range("M20").Value = 1.043132289345764
$ws.Range("N20").Value = 1.014061864220659
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.022756577709474
$ws.Range("D21").Value = 1.032693473329439
$ws.Range("E21").Value = 1.026617252812391
$ws.Range("F21").Value = 1.03867441567462
$ws.Range("I21").Value = 1.03206737848036
$ws.Range("J21").Value = 1.029439179347617
$ws.Range("K21").Value = 1.036285960653795
$ws.Range("L21").Value = 1.030232663387587
$ws.Range("M21").Value = 1.042244629038066
$ws.Range("N21").Value = 1.013825663948389
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.022057705860595
$ws.Range("D22").Value = 1.032134292878899
$ws.Range("E22").Value = 1.025958046980136
$ws.Range("F22").Value = 1.037984673980646
$ws.Range("I22").Value = 1.031917450943336
$ws.Range("J22").Value = 1.028993510700951
$ws.Range("K22").Value = 1.035858413958891
$ws.Range("L22").Value = 1.029706397883349
$ws.Range("M22").Value = 1.041686138790426
$ws.Range("N22").Value = 1.01367696724535
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.022428173340824
$ws.Range("D23").Value = 1.032430751452511
$ws.Range("E23").Value = 1.026307449178595
$ws.Range("F23").Value = 1.038350292913716
$ws.Range("I23").Value = 1.031997118987602
$ws.Range("J23").Value = 1.029229807899362
$ws.Range("K23").Value = 1.036085157691119
$ws.Range("L23").Value = 1.029985388989832
$ws.Range("M23").Value = 1.041982241091697
$ws.Range("N23").Value = 1.013755811945911
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.023887060570216
$ws.Range("D24").Value = 1.033597289555437
$ws.Range("E24").Value = 1.027684234844623
$ws.Range("F24").Value = 1.039790276826569
$ws.Range("I24").Value = 1.03230657895731
$ws.Range("J24").Value = 1.030159201921656
$ws.Range("K24").Value = 1.036975757108786
$ws.Range("L24").Value = 1.031083574356544
$ws.Range("M24").Value = 1.043147164541732
$ws.Range("N24").Value = 1.01406582094218
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.02558081747017
$ws.Range("D25").Value = 1.034949730664321
$ws.Range("E25").Value = 1.029284420068031
$ws.Range("F25").Value = 1.041462461078708
$ws.Range("I25").Value = 1.032657001359037
$ws.Range("J25").Value = 1.03123585609466
$ws.Range("K25").Value = 1.038004934502534
$ws.Range("L25").Value = 1.034110715006926
$ws.Range("M25").Value = 1.044497300547631
$ws.Range("N25").Value = 1.014424741867417
